$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# n_balaos / acertos_pedra counts updated
$ws.Range("A2").Value = 13
$ws.Range("B2").Value = 17
$ws.Range("C2").Value = 13

# precisao text percentage updated, keep it as plain text (not a numeric percent)
$ws.Range("D2").Value = "'76.5%"
$ws.Range("D2").Style = "Normal"
